$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stays text, since many values are numeric-looking
# (e.g. "211.85") and would otherwise be auto-converted to numbers by Excel,
# while others use a dotted-thousands format (e.g. "28.497.01") that must remain text.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.497.01"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "1.570.87"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "211.85"
$ws.Range("E5").Value = "  -0.62%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "46.09"
$ws.Range("E8").Value = "  +5.69%  "
$ws.Range("D9").Value = "24.14"
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "0.0881"
$ws.Range("E12").Value = "  -0.36%  "
$ws.Range("D13").Value = "1.796.89"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "1.575.75"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "0.520"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").Value = "3.68"
$ws.Range("E16").Value = "  -2.00%  "
$ws.Range("D17").Value = "28.496.38"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").Value = "227.01"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  -5.85%  "
$ws.Range("D24").Value = "9.10"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  +4.68%  "
$ws.Range("D26").Value = "151.17"
$ws.Range("E26").Value = "  -0.67%  "
$ws.Range("D27").Value = "14.95"
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("D28").Value = "6.44"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.0464"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.10"
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "1.391.11"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("E36").Value = "  -3.36%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("E39").Value = "  +3.72%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "0.531"
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("D43").Value = "0.792"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "0.979"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "63.01"
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("D48").Value = "1.708.80"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").Value = "85.99"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").Value = "0.0519"
$ws.Range("E51").Value = "  -1.27%  "
